$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = "ทดสอบระบบ"
$ws.Range("B8").Value = "ภาษาไทยพิมพ์ได้ป่าว"
$ws.Range("C8").Value = 50
$ws.Range("D8").Value = "Thai text"

$ws.Range("A8:D8").HorizontalAlignment = -4108
$ws.Range("A8:D8").VerticalAlignment = -4108
$ws.Rows.Item(8).RowHeight = 72

$ws.Range("F7").Select()
